# Update the Contact First Name And Last Name for row 3 (license 34-23421)
# from "Dakota Myers" to "Franz Ferdinand".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Franz Ferdinand"

# Move / record the active selection on the sheet, matching the saved view state.
[void]$ws.Range("D3").Select()
